# Add a new localization entry "strChkWindowPosition" to the Romanian (ro-RO)
# translation workbook. The sheet contains a table ("Tabla13") over
# B2:F192 that is kept sorted alphabetically by the "Key" column (C). The new
# key sorts in between the existing "strChkPower" (row 33) and
# "strDifferentiationAlgorithms" (row 34) entries, so it lands on worksheet
# row 34 and every row below it shifts down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at row 34 - this shifts rows 34..192 down to 35..193,
# carrying along their values, row heights and styles untouched.
$ws.Rows.Item(34).Insert()

# Grow the table (and its autofilter) so the new row becomes part of it.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("B2:F193"))

# Populate the new row with the new key/comment/english text; the Romanian
# translation column (F) is intentionally left blank, same as other
# not-yet-translated rows.
$ws.Range("B34").Value = "localization\strings"
$ws.Range("C34").Value = "strChkWindowPosition"
$ws.Range("D34").Value = 'In "settings" form, tab "User interface"'
$ws.Range("E34").Value = "Remember window position and size on startup"

# The pre-existing "strChkDlgPath" row (row 29, unaffected by the shift
# above since it is before the insertion point) also gets a comment added.
$ws.Range("D29").Value = 'In "settings" form, tab "User interface"'
